$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "kategori_nama" column (F) next to the existing data
$ws.Range("F1").Value = "kategori_nama"
$ws.Range("F2:F3").Value = "Personal Care"
$ws.Range("F4:F6").Value = "Baby Product"

# Match the header formatting used by the other headers (bold, centered)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Size the new column to fit its content
$ws.Columns.Item(6).ColumnWidth = 13.3

# Leave the selection where the user last clicked
$ws.Range("H6").Select()
